# Insert a new row at position 444, shifting existing rows 444+ down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(444).Insert()

# Populate the newly inserted row 444 with its data.
$ws.Range("A444").Value = 10
$ws.Range("B444").Value = "Vega Modelo de Temuco"
$ws.Range("C444").Value = "La Araucanía"
$ws.Range("D444").Value = 45218
$ws.Range("E444").Value = 9
$ws.Range("F444").Value = "Fruta"
$ws.Range("G444").Value = 100102
$ws.Range("H444").Value = "Cítricos"
$ws.Range("I444").Value = 100102006
$ws.Range("J444").Value = "Pomelo"
$ws.Range("K444").Value = "Start Ruby"
$ws.Range("L444").Value = "Primera"
$ws.Range("M444").Value = 220
$ws.Range("N444").Value = 14000
$ws.Range("O444").Value = 15000
$ws.Range("P444").Value = 14455
$ws.Range("Q444").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R444").Value = "Región de O'Higgins"
$ws.Range("S444").Value = 964
$ws.Range("T444").Value = 15
